$d = $word.ActiveDocument

$newTag = "{&quot;BindingKey&quot;:&quot;RefGroups2[0]&quot;,&quot;BindingType&quot;:&quot;Bookmark&quot;}"

for ($i = 1; $i -le $d.ContentControls.Count; $i++) {
    $cc = $d.ContentControls.Item($i)
    $cc.Tag = $newTag
}
